# v2p14. Compatible with MF-Swift v2212, updated hardpoints.
# Updates the nonlinear damping hardpoint table (F7:H8) on both
# "Sedan_HambaLG_f" and "Sedan_HambaLG_r" sheets, refreshes tab colors,
# and makes "Sedan_HambaLG_f" the active/selected sheet.

$wb = $excel.ActiveWorkbook

$wsF = $wb.Worksheets.Item("Sedan_HambaLG_f")
$wsR = $wb.Worksheets.Item("Sedan_HambaLG_r")

# ---------------------------------------------------------------------
# Sedan_HambaLG_f (sheet1): damping hardpoints
# ---------------------------------------------------------------------
$wsF.Range("F7").NumberFormat = "0.000"
$wsF.Range("F7").Value = -0.002655714285714287
$wsF.Range("G7").Value = 0.62
$wsF.Range("H7").Value = 0.65

$wsF.Range("F8").NumberFormat = "0.000"
$wsF.Range("F8").Value = 0.05516642857142858
$wsF.Range("G8").Value = 0.85
$wsF.Range("H8").Value = 0.19

# ---------------------------------------------------------------------
# Sedan_HambaLG_r (sheet2): damping hardpoints
# ---------------------------------------------------------------------
$wsR.Range("F7").Value = 0.002655714285714287
$wsR.Range("G7").Value = 0.62
$wsR.Range("H7").Value = 0.65

$wsR.Range("F8").Value = -0.05516642857142858
$wsR.Range("G8").Value = 0.85
$wsR.Range("H8").Value = 0.19

# ---------------------------------------------------------------------
# Tab colors: Gold Accent4 Lighter80% -> Blue Accent5 Darker25%
# (theme color write-back isn't fully supported by this runtime, so the
# resolved RGB equivalent is applied instead for visual fidelity)
# ---------------------------------------------------------------------
$wsF.Tab.Color = 11957550
$wsR.Tab.Color = 11957550

# ---------------------------------------------------------------------
# Selection / active sheet: Sedan_HambaLG_f becomes the active tab with
# F7:H8 selected; Sedan_HambaLG_r keeps a plain single-cell selection.
# ---------------------------------------------------------------------
$wsR.Activate()
$wsR.Range("E15").Select()

$wsF.Activate()
$wsF.Range("F7:H8").Select()
